# Update cfb_weather.xlsx with Timestamp 2025-12-01T05:15:41.461386
# This script updates the FBS sheet (sheet1) and Other sheet (sheet2)
# with refreshed weather / odds data, matching a later data-pull run.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("FBS")
$ws2 = $wb.Worksheets.Item("Other")

# ---- FBS sheet (sheet1) ----
$ws1.Range('AA2').Value = [double]-3.5
$ws1.Range('AE2').Value = [double]0
$ws1.Range('AF2').Value = [double]0
$ws1.Range('AK2').Value = '2025-12-01T05:15:41.461386'
$ws1.Range('M2').Value = 'E'
$ws1.Range('N2').Value = 'E'
$ws1.Range('O2').Value = [double]50.3
$ws1.Range('P2').Value = [double]16
$ws1.Range('Q2').Value = 'E'
$ws1.Range('S2').Value = [double]-3.5
$ws1.Range('U2').Value = [double]9.199999999999999
$ws1.Range('V2').Value = '43.6028839, -116.1958882'
$ws1.Range('W2').Value = [double]57.5
$ws1.Range('AA3').Value = [double]-22
$ws1.Range('AF3').Value = [double]0
$ws1.Range('AK3').Value = '2025-12-01T05:15:41.461386'
$ws1.Range('K3').Value = [double]6.2
$ws1.Range('O3').Value = [double]28.4
$ws1.Range('P3').Value = [double]7.4
$ws1.Range('R3').Value = [double]0
$ws1.Range('S3').Value = [double]-0.2
$ws1.Range('T3').Value = [double]-0.45
$ws1.Range('U3').Value = [double]1.2
$ws1.Range('V3').Value = '38.4352919, -78.8729349'
$ws1.Range('Z3').Value = [double]-110
$ws1.Range('AB4').Value = [double]-1
$ws1.Range('AF4').Value = [double]-0.5
$ws1.Range('AK4').Value = '2025-12-01T05:15:41.461386'
$ws1.Range('K4').Value = [double]4.8
$ws1.Range('M4').Value = 'SSW'
$ws1.Range('N4').Value = 'SW'
$ws1.Range('O4').Value = [double]47
$ws1.Range('P4').Value = [double]1.3
$ws1.Range('Q4').Value = 'SW'
$ws1.Range('R4').Value = [double]0
$ws1.Range('U4').Value = [double]-3.5
$ws1.Range('V4').Value = '33.8201052, -85.76647'
$ws1.Range('X4').Value = [double]-105
$ws1.Range('Z4').Value = [double]-110
$ws1.Range('AA5').Value = [double]2.5
$ws1.Range('AE5').Value = [double]0
$ws1.Range('AF5').Value = [double]0
$ws1.Range('AK5').Value = '2025-12-01T05:15:41.461386'
$ws1.Range('K5').Value = [double]10.6
$ws1.Range('M5').Value = 'S'
$ws1.Range('N5').Value = 'S'
$ws1.Range('O5').Value = [double]59.06
$ws1.Range('P5').Value = [double]7.6
$ws1.Range('Q5').Value = 'S'
$ws1.Range('R5').Value = [double]0.6000000000000001
$ws1.Range('U5').Value = [double]-3
$ws1.Range('V5').Value = '29.944616, -90.116692'
$ws1.Range('W5').Value = [double]67.5
$ws1.Range('X5').Value = [double]-115
$ws1.Range('Z5').Value = [double]-110
$ws1.Range('AE6').Value = [double]0
$ws1.Range('AK6').Value = '2025-12-01T05:15:41.461386'
$ws1.Range('K6').Value = [double]12
$ws1.Range('M6').Value = 'N'
$ws1.Range('N6').Value = 'N'
$ws1.Range('O6').Value = [double]31.34
$ws1.Range('P6').Value = [double]9.699999999999999
$ws1.Range('Q6').Value = 'N'
$ws1.Range('U6').Value = [double]-2.3
$ws1.Range('V6').Value = '42.2860064, -85.6007573'
$ws1.Range('W6').Value = [double]43.5
$ws1.Range('X6').Value = [double]-105
$ws1.Range('AA7').Value = [double]-3
$ws1.Range('AF7').Value = [double]0
$ws1.Range('AK7').Value = '2025-12-01T05:15:41.461386'
$ws1.Range('K7').Value = [double]4.6
$ws1.Range('M7').Value = 'W'
$ws1.Range('N7').Value = 'WSW'
$ws1.Range('O7').Value = [double]36.02
$ws1.Range('P7').Value = [double]1.3
$ws1.Range('Q7').Value = 'SW'
$ws1.Range('V7').Value = '38.0311801, -78.5137897'
$ws1.Range('X7').Value = [double]-115

# ---- Other sheet (sheet2) ----
$ws2.Range('B1').Value = 'Home Team'
$ws2.Range('C1').Value = 'Away Team'
$ws2.Range('A2').Value = 'North Dakota vs Tarleton State'
$ws2.Range('B2').Value = 'Tarleton State'
$ws2.Range('C2').Value = 'North Dakota'
$ws2.Range('O2').Value = 'NE'
$ws2.Range('P2').Value = 'NE'
$ws2.Range('Q2').Value = [double]70.88000000000001
$ws2.Range('R2').Value = [double]25.1
$ws2.Range('S2').Value = 'NE'
$ws2.Range('X2').Value = '32.2191836, -98.2130634'
$ws2.Range('A3').Value = 'Yale vs Montana State'
$ws2.Range('B3').Value = 'Montana State'
$ws2.Range('C3').Value = 'Yale'
$ws2.Range('J3').Value = [double]1502.206045159
$ws2.Range('K3').Value = [double]42.68
$ws2.Range('L3').Value = [double]53.64
$ws2.Range('N3').Value = [double]1973
$ws2.Range('O3').Value = 'ESE'
$ws2.Range('P3').Value = 'E'
$ws2.Range('Q3').Value = [double]19.34
$ws2.Range('R3').Value = [double]1
$ws2.Range('S3').Value = 'E'
$ws2.Range('T3').Value = [double]0
$ws2.Range('U3').Value = [double]-1.33
$ws2.Range('V3').Value = [double]-3.5
$ws2.Range('X3').Value = '45.659048, -111.049547'
$ws2.Range('A4').Value = 'Villanova vs Lehigh'
$ws2.Range('B4').Value = 'Lehigh'
$ws2.Range('C4').Value = 'Villanova'
$ws2.Range('J4').Value = [double]-37.06062315000001
$ws2.Range('K4').Value = [double]54.29
$ws2.Range('L4').Value = [double]55.05
$ws2.Range('N4').Value = [double]1988
$ws2.Range('O4').Value = 'ENE'
$ws2.Range('P4').Value = 'ENE'
$ws2.Range('Q4').Value = [double]34.88
$ws2.Range('R4').Value = [double]5.2
$ws2.Range('S4').Value = 'ENE'
$ws2.Range('U4').Value = [double]0
$ws2.Range('V4').Value = [double]0
$ws2.Range('X4').Value = '40.5890837, -75.3553874'
$ws2.Range('A5').Value = 'South Dakota vs Mercer'
$ws2.Range('B5').Value = 'Mercer'
$ws2.Range('C5').Value = 'South Dakota'
$ws2.Range('O5').Value = 'NNW'
$ws2.Range('P5').Value = 'NNW'
$ws2.Range('Q5').Value = [double]59.54
$ws2.Range('R5').Value = [double]2.8
$ws2.Range('S5').Value = 'NNW'
$ws2.Range('X5').Value = '32.8262075, -83.6522485'
$ws2.Range('A6').Value = 'Abilene Christian vs Stephen F. Austin'
$ws2.Range('B6').Value = 'Stephen F. Austin'
$ws2.Range('C6').Value = 'Abilene Christian'
$ws2.Range('O6').Value = 'NNE'
$ws2.Range('P6').Value = 'NNE'
$ws2.Range('Q6').Value = [double]66.26000000000001
$ws2.Range('R6').Value = [double]11.3
$ws2.Range('S6').Value = 'NNE'
$ws2.Range('X6').Value = '31.625719, -94.6444034'
$ws2.Range('A7').Value = 'South Dakota State vs Montana'
$ws2.Range('B7').Value = 'Montana'
$ws2.Range('C7').Value = 'South Dakota State'
$ws2.Range('O7').Value = 'WNW'
$ws2.Range('P7').Value = 'WNW'
$ws2.Range('Q7').Value = [double]31.1
$ws2.Range('R7').Value = [double]4.5
$ws2.Range('S7').Value = 'WNW'
$ws2.Range('T7').Value = [double]0
$ws2.Range('X7').Value = '46.8638753, -113.9815042'
$ws2.Range('A8').Value = 'Rhode Island vs UC Davis'
$ws2.Range('B8').Value = 'UC Davis'
$ws2.Range('C8').Value = 'Rhode Island'
$ws2.Range('O8').Value = 'NE'
$ws2.Range('P8').Value = 'E'
$ws2.Range('Q8').Value = [double]49.63999999999999
$ws2.Range('R8').Value = [double]1.7
$ws2.Range('S8').Value = 'NE'
$ws2.Range('X8').Value = '38.5365266, -121.7627936'
